$wb = $excel.ActiveWorkbook

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/75be731a63683d49190e51e40793af59acf1723b/e2e/a.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/75be731a63683d49190e51e40793af59acf1723b/e2e/b.md"

# ---- Overview sheet: status text changes from "Ready for handoff" to
# "Handed back: in sync with en-US" for both zh-cn (E) and de-de (F) columns.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

# ---- zh-cn sheet: the handback for a.md finished -- fill in the target
# file, handback file and handback datetime, same status text change,
# and add hyperlinks on the Latest Target File cells.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("I2").Value = "a.md"
$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-19 04:32:34"
$wsZh.Range("I3").Value = "a.md"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-19 04:32:34"
$wsZh.Columns.Item(3).ColumnWidth = 29.14
$wsZh.Columns.Item(10).ColumnWidth = 40

# Rebuild the hyperlinks collection so relationship ids come out in the
# same row-major order as the handoff hyperlinks (A2, I2, A3, I3).
$wsZh.Hyperlinks.Delete()
$null = $wsZh.Hyperlinks.Add($wsZh.Range("A2"), $urlA, [System.Type]::Missing, [System.Type]::Missing, "a.md")
$null = $wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlA, [System.Type]::Missing, [System.Type]::Missing, "a.md")
$null = $wsZh.Hyperlinks.Add($wsZh.Range("A3"), $urlB, [System.Type]::Missing, [System.Type]::Missing, "b.md")
$null = $wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlA, [System.Type]::Missing, [System.Type]::Missing, "a.md")

# ---- de-de sheet: same shape, different handback file name / datetime.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("I2").Value = "a.md"
$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-19 04:32:41"
$wsDe.Range("I3").Value = "a.md"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-19 04:32:41"
$wsDe.Columns.Item(3).ColumnWidth = 29.14
$wsDe.Columns.Item(10).ColumnWidth = 40

$wsDe.Hyperlinks.Delete()
$null = $wsDe.Hyperlinks.Add($wsDe.Range("A2"), $urlA, [System.Type]::Missing, [System.Type]::Missing, "a.md")
$null = $wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlA, [System.Type]::Missing, [System.Type]::Missing, "a.md")
$null = $wsDe.Hyperlinks.Add($wsDe.Range("A3"), $urlB, [System.Type]::Missing, [System.Type]::Missing, "b.md")
$null = $wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlA, [System.Type]::Missing, [System.Type]::Missing, "a.md")
